$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.247.63'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '3.502.89'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('D7').Value = '3.503.74'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -0.25%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.375'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.73%  '
$ws.Range('D13').Value = '4.098.46'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').Value = '3.501.44'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.48%  '
$ws.Range('D18').Value = '64.266.39'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('E22').Value = '  -2.66%  '
$ws.Range('E23').Value = '  -1.98%  '
$ws.Range('D24').Value = '3.640.99'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('E28').Value = '  +2.94%  '
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('D34').Value = '3.522.34'
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.146'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.87'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('E40').Value = '  -4.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '164.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.00%  '
$ws.Range('E42').Value = '  -4.36%  '
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').Value = '2.470.27'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.918'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.74%  '
